$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 60, pushing existing rows 60-72 down to 61-73.
$ws.Rows.Item(60).Insert()

# Fill in the new weekly record in row 60.
$ws.Cells.Item(60, 1).Value2 = 7
$ws.Cells.Item(60, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(60, 3).Value2 = "Ñuble"
$ws.Cells.Item(60, 4).Value2 = 44617
$ws.Cells.Item(60, 5).Value2 = 16
$ws.Cells.Item(60, 6).Value2 = 100112031
$ws.Cells.Item(60, 7).Value2 = "Poroto verde"
$ws.Cells.Item(60, 8).Value2 = "Sin especificar"
$ws.Cells.Item(60, 9).Value2 = "Primera"
$ws.Cells.Item(60, 10).Value2 = 120
$ws.Cells.Item(60, 11).Value2 = 28000
$ws.Cells.Item(60, 12).Value2 = 29000
$ws.Cells.Item(60, 13).Value2 = 28500
$ws.Cells.Item(60, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(60, 15).Value2 = "Región del Maule"
$ws.Cells.Item(60, 16).Value2 = 1140
$ws.Cells.Item(60, 17).Value2 = 25
$ws.Cells.Item(60, 18).Value2 = "Hortaliza"
